$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.823.57'
$ws.Range('E2').Value = '  +2.53%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.159.44'
$ws.Range('E3').Value = '  +2.73%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '537.67'
$ws.Range('E5').Value = '  +2.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.44'
$ws.Range('E6').Value = '  +2.69%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.512'
$ws.Range('E8').Value = '  +8.55%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.35'
$ws.Range('E9').Value = '  +1.65%  '
$ws.Range('E10').Value = '  +3.75%  '
$ws.Range('E11').Value = '  +5.06%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.140'
$ws.Range('E12').Value = '  +1.98%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.701.74'
$ws.Range('E13').Value = '  +2.77%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.94'
$ws.Range('E14').Value = '  +3.08%  '
$ws.Range('E15').Value = '  +6.51%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '58.832.83'
$ws.Range('E16').Value = '  +2.54%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.169.16'
$ws.Range('E17').Value = '  +3.27%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.23'
$ws.Range('E18').Value = '  +5.99%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.04'
$ws.Range('E19').Value = '  +4.88%  '
$ws.Range('E20').Value = '  +5.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '373.18'
$ws.Range('E21').Value = '  +6.81%  '
$ws.Range('E22').Value = '  +1.96%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.998'
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.95'
$ws.Range('E24').Value = '  +2.75%  '
$ws.Range('E25').Value = '  +3.09%  '
$ws.Range('E26').Value = '  +1.55%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.09'
$ws.Range('E28').Value = '  +14.14%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0879'
$ws.Range('E29').Value = '  +3.94%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.90'
$ws.Range('E30').Value = '  +2.54%  '
$ws.Range('B31').Value = 'RenderToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.17'
$ws.Range('E31').Value = '  +4.44%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.97'
$ws.Range('E32').Value = '  +4.49%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.21'
$ws.Range('E33').Value = '  +8.10%  '
$ws.Range('E34').Value = '  +4.11%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '159.91'
$ws.Range('E35').Value = '  +0.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.25'
$ws.Range('E36').Value = '  +4.72%  '
$ws.Range('E37').Value = '  +10.86%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '25.33'
$ws.Range('E38').Value = '  -0.41%  '
$ws.Range('E39').Value = '  +7.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.657.01'
$ws.Range('E40').Value = '  +10.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0686'
$ws.Range('E41').Value = '  +4.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.19'
$ws.Range('E42').Value = '  +4.33%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '38.85'
$ws.Range('E43').Value = '  +5.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.709'
$ws.Range('E44').Value = '  +2.99%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0284'
$ws.Range('E45').Value = '  +9.15%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.00'
$ws.Range('E46').Value = '  -0.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.199.65'
$ws.Range('E47').Value = '  +2.77%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.103'
$ws.Range('E48').Value = '  +12.19%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.990'
$ws.Range('E49').Value = '  +4.75%  '
$ws.Range('E50').Value = '  +4.04%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '20.33'
$ws.Range('E51').Value = '  +5.36%  '
